$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("grandes regiões e unidades da federação") is an empty header row with no
# numeric data. The fix removes this row entirely: the region labels/data that used
# to live in rows 7-38 shift up to rows 6-37, and the sheet shrinks from A1:G38 to A1:G37.
$ws.Rows.Item(6).Delete()
